$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column F (dSF) values per repulled data / mean calculation
$ws.Range("F2").Value = -7
$ws.Range("F3").Value = 4
$ws.Range("F4").Value = -2
$ws.Range("F6").Value = -5
$ws.Range("F9").Value = -6
$ws.Range("F11").Value = -5
$ws.Range("F12").Value = -4
$ws.Range("F15").Value = -1
